# This script applies timetable corrections across the mon/tue/wed/thur/fri
# sheets of the final_timetable workbook, per the authored diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: mon
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mon")

$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("I4").Value = "CSC423"

$ws.Range("F7").Value = "CSC111"
$ws.Range("G7").Value = "CSC111"

$ws.Range("J15").Value = "BIO111"

$ws.Range("D17").Value = "CSC424"

$ws.Range("K18").Value = ""

$ws.Range("H19").Value = "CSC424"
$ws.Range("I19").Value = "CSC424"

$ws.Range("E20").Value = "MAT111"
$ws.Range("F20").Value = "MAT111"

$ws.Range("F21").Value = "GST111"
$ws.Range("G21").Value = "GST111"
$ws.Range("H21").Value = ""
$ws.Range("I21").Value = ""
$ws.Range("J21").Value = ""

# ---------------------------------------------------------------------
# Sheet: tue
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("tue")

$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("J15").Value = ""

$ws.Range("H17").Value = ""

$ws.Range("K20").Value = "MAT111"

$ws.Range("D21").Value = "MAT112"
$ws.Range("E21").Value = "MAT112"
$ws.Range("I21").Value = "CST111"
$ws.Range("J21").Value = "CST111"

# ---------------------------------------------------------------------
# Sheet: wed
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("wed")

$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""

$ws.Range("I7").Value = "CSC425"

$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""

$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""

$ws.Range("H13").Value = ""
$ws.Range("I13").Value = ""

$ws.Range("D14").Value = "BIO111"
$ws.Range("E14").Value = "BIO111"

$ws.Range("D20").Value = "CIT111"
$ws.Range("E20").Value = "CIT111"
$ws.Range("H20").Value = ""
$ws.Range("I20").Value = "EDS421"
$ws.Range("J20").Value = ""
$ws.Range("K20").Value = ""

$ws.Range("D21").Value = ""

# ---------------------------------------------------------------------
# Sheet: thur
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("thur")

$ws.Range("G2").Value = "CSC425"
$ws.Range("H2").Value = "CSC425"

$ws.Range("E13").Value = "CSC424"

$ws.Range("D17").Value = "CSC424"
$ws.Range("E17").Value = "CSC424"
$ws.Range("F17").Value = "CSC442"
$ws.Range("G17").Value = ""
$ws.Range("H17").Value = ""

$ws.Range("K18").Value = ""

$ws.Range("J20").Value = ""
$ws.Range("K20").Value = ""

$ws.Range("F21").Value = ""
$ws.Range("G21").Value = ""

$ws.Range("E24").Value = "CSC423"
$ws.Range("F24").Value = "CSC423"

# ---------------------------------------------------------------------
# Sheet: fri
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("fri")

$ws.Range("F7").Value = "CIS421"
$ws.Range("G7").Value = "CIS421"

$ws.Range("C28").Value = "CSC441"
$ws.Range("D28").Value = "CSC441"

Write-Output "timetable updates applied"
